$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8 - this shifts existing rows 8-13 down to 9-14
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(8, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(8, 4).Value = 44484
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = 100112026
$ws.Cells.Item(8, 7).Value = 'Haba'
$ws.Cells.Item(8, 8).Value = 'Sin especificar'
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 900
$ws.Cells.Item(8, 11).Value = 750
$ws.Cells.Item(8, 12).Value = 800
$ws.Cells.Item(8, 13).Value = 775
$ws.Cells.Item(8, 14).Value = '$/kilo'
$ws.Cells.Item(8, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(8, 16).Value = 775
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = 'Hortaliza'

# Match the date-style formatting used by column D in the other rows
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
